$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.908.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.611.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0755"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.912.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.141.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.798"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.770.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.415"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.26%  "
